# Updates cryptos list values (prices and 1h volume deltas) to match the
# Thu Sep 26 11:53:09 UTC 2024 GitHub Actions refresh of cryptos.xlsx.
# Row 10/11 and 38/39 swap coin identity (Toncoin<->Cardano, EthereumClassic<->Stacks).
#
# Price cells (column D) that look like plain decimal numbers are written with a
# leading apostrophe so Excel stores them as text (quotePrefix) instead of silently
# converting them to floating point numbers and losing the exact formatting
# (e.g. trailing zeros such as "1.00" or "169.00").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.398.63"
$ws.Range("E2").Value = "  +1.02%  "

$ws.Range("D3").Value = "2.629.89"
$ws.Range("E3").Value = "  +0.22%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'596.07"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").Value = "'152.80"
$ws.Range("E6").Value = "  +1.37%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +3.68%  "

$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "'0.398"
$ws.Range("E10").Value = "  +4.13%  "

$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "'5.84"
$ws.Range("E11").Value = "  +2.77%  "

$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("D14").Value = "3.100.23"
$ws.Range("E14").Value = "  +0.26%  "

$ws.Range("D15").Value = "64.290.63"
$ws.Range("E15").Value = "  +1.15%  "

$ws.Range("D16").Value = "'0.0000170"
$ws.Range("E16").Value = "  +10.53%  "

$ws.Range("D17").Value = "2.623.10"
$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("E18").Value = "  -0.12%  "

$ws.Range("D19").Value = "'4.77"
$ws.Range("E19").Value = "  +2.23%  "

$ws.Range("D20").Value = "'349.59"
$ws.Range("E20").Value = "  +0.50%  "

$ws.Range("E21").Value = "  +3.31%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").Value = "'67.69"
$ws.Range("E23").Value = "  +1.92%  "

$ws.Range("D24").Value = "'1.72"
$ws.Range("E24").Value = "  -0.51%  "

$ws.Range("E25").Value = "  +0.21%  "

$ws.Range("E26").Value = "  -0.52%  "

$ws.Range("D27").Value = "'8.31"
$ws.Range("E27").Value = "  +0.70%  "

$ws.Range("D28").Value = "'547.78"
$ws.Range("E28").Value = "  -3.52%  "

$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").Value = "'0.995"
$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("D31").Value = "0.0₃0915"
$ws.Range("E31").Value = "  +7.90%  "

$ws.Range("E32").Value = "  +1.31%  "

$ws.Range("E33").Value = "  +4.10%  "

$ws.Range("D34").Value = "'5.55"
$ws.Range("E34").Value = "  +6.06%  "

$ws.Range("E35").Value = "  +1.36%  "

$ws.Range("E36").Value = "  +3.22%  "

$ws.Range("D37").Value = "'165.23"
$ws.Range("E37").Value = "  -2.19%  "

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'2.01"
$ws.Range("E38").Value = "  +3.21%  "

$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "'20.07"
$ws.Range("E39").Value = "  +3.35%  "

$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").Value = "'169.00"
$ws.Range("E42").Value = "  +1.42%  "

$ws.Range("D43").Value = "'41.77"
$ws.Range("E43").Value = "  +4.61%  "

$ws.Range("E44").Value = "  +5.02%  "

$ws.Range("D45").Value = "'23.37"
$ws.Range("E45").Value = "  +8.02%  "

$ws.Range("D46").Value = "'2.24"
$ws.Range("E46").Value = "  +12.18%  "

$ws.Range("E47").Value = "  -0.83%  "

$ws.Range("D48").Value = "'0.640"
$ws.Range("E48").Value = "  +1.42%  "

$ws.Range("E49").Value = "  +0.84%  "

$ws.Range("D50").Value = "'0.0979"
$ws.Range("E50").Value = "  +1.41%  "

$ws.Range("D51").Value = "'19.38"
$ws.Range("E51").Value = "  -0.31%  "
